$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.621.47"
$ws.Cells.Item(2, 5).Value = "  -3.64%  "
$ws.Cells.Item(3, 4).Value = "3.332.71"
$ws.Cells.Item(3, 5).Value = "  -4.28%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$c = $ws.Cells.Item(5, 4)
$c.Value = "'181.46"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -8.32%  "
$c = $ws.Cells.Item(6, 4)
$c.Value = "'532.22"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.77%  "
$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.608"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.74%  "
$ws.Cells.Item(8, 4).Value = "3.327.65"
$ws.Cells.Item(8, 5).Value = "  -4.26%  "
$ws.Cells.Item(9, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.615"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -5.11%  "
$c = $ws.Cells.Item(11, 4)
$c.Value = "'58.67"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -5.98%  "
$ws.Cells.Item(12, 5).Value = "  -5.29%  "
$ws.Cells.Item(13, 5).Value = "  -2.20%  "
$c = $ws.Cells.Item(14, 4)
$c.Value = "'9.18"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -5.85%  "
$ws.Cells.Item(15, 4).Value = "3.867.86"
$ws.Cells.Item(15, 5).Value = "  -4.60%  "
$ws.Cells.Item(16, 4).Value = "3.335.61"
$ws.Cells.Item(16, 5).Value = "  -4.52%  "
$ws.Cells.Item(17, 5).Value = "  -4.40%  "
$ws.Cells.Item(18, 4).Value = "64.594.28"
$ws.Cells.Item(18, 5).Value = "  -3.33%  "
$c = $ws.Cells.Item(19, 4)
$c.Value = "'17.60"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.30%  "
$c = $ws.Cells.Item(20, 4)
$c.Value = "'11.22"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -4.04%  "
$c = $ws.Cells.Item(21, 4)
$c.Value = "'0.967"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -4.76%  "
$c = $ws.Cells.Item(22, 4)
$c.Value = "'377.11"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.73%  "
$c = $ws.Cells.Item(23, 4)
$c.Value = "'3.83"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.24%  "
$ws.Cells.Item(24, 5).Value = "  -6.76%  "
$c = $ws.Cells.Item(25, 4)
$c.Value = "'81.28"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.14%  "
$ws.Cells.Item(26, 5).Value = "  +2.79%  "
$ws.Cells.Item(27, 2).Value = "LEO"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Cells.Item(27, 4)
$c.Value = "'6.09"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.16%  "
$ws.Cells.Item(28, 2).Value = "ImmutableX"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Cells.Item(28, 4)
$c.Value = "'2.70"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.07%  "
$ws.Cells.Item(29, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Cells.Item(29, 4)
$c.Value = "'11.52"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -4.56%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Cells.Item(30, 4)
$c.Value = "'8.44"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -3.45%  "
$ws.Cells.Item(31, 2).Value = "EthereumClassic"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Cells.Item(31, 4)
$c.Value = "'29.20"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -5.09%  "
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Cells.Item(32, 4)
$c.Value = "'653.46"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -3.02%  "
$ws.Cells.Item(33, 2).Value = "NEARProtocol"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(33, 4)
$c.Value = "'6.71"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.69%  "
$ws.Cells.Item(34, 2).Value = "Cosmos"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(34, 4)
$c.Value = "'11.35"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.37%  "
$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.107"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -3.53%  "
$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Cells.Item(36, 4)
$c.Value = "'59.73"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -6.07%  "
$ws.Cells.Item(37, 2).Value = "Dai"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Cells.Item(37, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.02%  "
$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.394"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.56%  "
$ws.Cells.Item(39, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Cells.Item(39, 4)
$c.Value = "'37.10"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.34%  "
$ws.Cells.Item(40, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(40, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.02%  "
$ws.Cells.Item(41, 2).Value = "PEPE"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(41, 4).Value = "0.0₃0713"
$ws.Cells.Item(41, 5).Value = "  +7.38%  "
$ws.Cells.Item(42, 2).Value = "Kaspa"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.126"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.99%  "
$ws.Cells.Item(43, 2).Value = "Maker"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(43, 4).Value = "2.946.21"
$ws.Cells.Item(43, 5).Value = "  -2.91%  "
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'2.52"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.17%  "
$ws.Cells.Item(45, 2).Value = "ThetaToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'2.73"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -7.39%  "
$ws.Cells.Item(46, 2).Value = "VeChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.0402"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.60%  "
$ws.Cells.Item(47, 2).Value = "WEMIXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Cells.Item(47, 4)
$c.Value = "'2.67"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -3.51%  "
$ws.Cells.Item(48, 2).Value = "Stacks"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'2.81"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +8.00%  "
$ws.Cells.Item(49, 2).Value = "ApeXProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Cells.Item(49, 4)
$c.Value = "'3.03"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +5.12%  "
$ws.Cells.Item(50, 2).Value = "Stellar"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(50, 4)
$c.Value = "'0.127"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.75%  "
$ws.Cells.Item(51, 2).Value = "dogwifhat"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(51, 4)
$c.Value = "'2.53"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -4.33%  "
